# Generate Report for Handback
# Update the timestamp strings recorded on the handback-status report.
#
# Mapping of changed shared strings -> cells:
#   "Overview"!G2 (Latest HO Xliff Generate Date)     2016-08-26 17:04:56 -> 2016-08-26 17:06:11
#       (shared string also used by "de-de"!H2, which updates alongside it)
#   "zh-cn"!H2   (Correspond Handoff Datetime)        2016-08-26 17:04:51 -> 2016-08-26 17:06:00
#   "zh-cn"!K2   (Correspond Handback DateTime)       2016-08-26 17:05:31 -> 2016-08-26 17:06:28
#   "de-de"!K2   (Correspond Handback DateTime)       2016-08-26 17:05:37 -> 2016-08-26 17:06:34

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-26 17:06:11"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-26 17:06:00"
$wsZhCn.Range("K2").Value = "2016-08-26 17:06:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
# de-de!H2 shares the same underlying text as Overview!G2, so it is
# refreshed to the same new value as well.
$wsDeDe.Range("H2").Value = "2016-08-26 17:06:11"
$wsDeDe.Range("K2").Value = "2016-08-26 17:06:34"
